$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 8 new audio records (rows 10-17) to the atlas_list sheet,
# mirroring the existing species_scientific_name / track_url / track_title /
# author_name / author_url / copyright_image / copyright_url layout.

$newRows = New-Object 'object[,]' 8,7

# Row 10 - Acanthiza pusilla (Noosa Heads)
$newRows[0,0] = "Acanthiza pusilla"
$newRows[0,1] = "https://www.xeno-canto.org/157488/download"
$newRows[0,2] = "Noosa Heads, Queensland, Australia"
$newRows[0,3] = "Fernand Deroussen"
$newRows[0,4] = "https://www.xeno-canto.org/contributor/UXGZWVYDFE"
$newRows[0,5] = "assets/misc/cc.png"
$newRows[0,6] = "https://creativecommons.org/licenses/by-nc-nd/3.0/"

# Row 11 - Acanthiza pusilla (Comerong Island)
$newRows[1,0] = "Acanthiza pusilla"
$newRows[1,1] = "https://www.xeno-canto.org/396904/download"
$newRows[1,2] = "Comerong Island, New South Wales, Australia"
$newRows[1,3] = "Greg McLachlan"
$newRows[1,4] = "https://www.xeno-canto.org/contributor/MXMFTGKZDR"
$newRows[1,5] = "assets/misc/cc.png"
$newRows[1,6] = "https://creativecommons.org/licenses/by-nc-sa/4.0/"

# Row 12 - Acanthiza reguloides (Moggill State Forest)
$newRows[2,0] = "Acanthiza reguloides"
$newRows[2,1] = "https://www.xeno-canto.org/85990/download"
$newRows[2,2] = "Moggill State Forest, Queensland, Australia"
$newRows[2,3] = "Nick Leseberg"
$newRows[2,4] = "https://www.xeno-canto.org/contributor/FWTRWUQQAA"
$newRows[2,5] = "assets/misc/cc.png"
$newRows[2,6] = "https://creativecommons.org/licenses/by-nc-nd/2.5/"

# Row 13 - Acanthiza reguloides (Nangar National Park)
$newRows[3,0] = "Acanthiza reguloides"
$newRows[3,1] = "https://www.xeno-canto.org/365188/download"
$newRows[3,2] = "Nangar National Park, New South Wales, Australia"
$newRows[3,3] = "Greg McLachlan"
$newRows[3,4] = "https://www.xeno-canto.org/contributor/MXMFTGKZDR"
$newRows[3,5] = "assets/misc/cc.png"
$newRows[3,6] = "https://creativecommons.org/licenses/by-nc-sa/4.0/"

# Row 14 - Acanthorhynchus tenuirostris (Abercrombie Caves)
$newRows[4,0] = "Acanthorhynchus tenuirostris"
$newRows[4,1] = "https://www.xeno-canto.org/357302/download"
$newRows[4,2] = "Abercrombie Caves, New South Wales, Australia"
$newRows[4,3] = "Greg McLachlan"
$newRows[4,4] = "https://www.xeno-canto.org/contributor/MXMFTGKZDR"
$newRows[4,5] = "assets/misc/cc.png"
$newRows[4,6] = "https://creativecommons.org/licenses/by-nc-sa/4.0/"

# Row 15 - Accipiter cirrocephalus (Royal National Park)
$newRows[5,0] = "Accipiter cirrocephalus"
$newRows[5,1] = "https://www.xeno-canto.org/209907/download"
$newRows[5,2] = "Royal National Park, New South Wales, Australia"
$newRows[5,3] = "Marc Anderson"
$newRows[5,4] = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$newRows[5,5] = "assets/misc/cc.png"
$newRows[5,6] = "https://creativecommons.org/licenses/by-nc-nd/4.0/"

# Row 16 - Accipiter fasciatus (Wollemi National Park)
$newRows[6,0] = "Accipiter fasciatus"
$newRows[6,1] = "https://www.xeno-canto.org/340225/download"
$newRows[6,2] = "Wollemi National Park, New South Wales, Australia"
$newRows[6,3] = "Marc Anderson"
$newRows[6,4] = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$newRows[6,5] = "assets/misc/cc.png"
$newRows[6,6] = "https://creativecommons.org/licenses/by-nc-nd/4.0/"

# Row 17 - Accipiter novaehollandiae (Eungella)
$newRows[7,0] = "Accipiter novaehollandiae"
$newRows[7,1] = "https://www.xeno-canto.org/434627/download"
$newRows[7,2] = "Eungella, Queensland, Australia"
$newRows[7,3] = "Marc Anderson"
$newRows[7,4] = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$newRows[7,5] = "assets/misc/cc.png"
$newRows[7,6] = "https://creativecommons.org/licenses/by-nc-nd/4.0/"

$ws.Range("A10:G17").Value = $newRows

# Re-fit column A now that a longer species name has been added.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Match the author's final selection/active cell.
$ws.Range("F17").Select() | Out-Null
